# Apply the text updates to KCOR_Deviation_by_Age_and_Frailty sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# K2: "The purpose of this table is to show the KCOR deviation due to frailty rate and age of the cohort."
#     -> "... over a 12 month period."
$ws.Range("K2").Value = "The purpose of this table is to show the KCOR deviation due to frailty rate and age of the cohort over a 12 month period."

# A1: "Given a frailty ratio and age, shows the KCOR distortion for that cohort"
#     -> "... over 12 months"
$ws.Range("A1").Value = "Given a frailty ratio and age, shows the KCOR distortion for that cohort over 12 months"

# K4 and K6 text content is unchanged (only shared-string ordering changed upstream,
# which is an internal storage detail, not a visible value change), but we re-set
# them explicitly to be safe / deterministic.
$ws.Range("K4").Value = "So an 85 year old cohort with a 3X frailty will have a -6.5% impact on slope making the vaccinated look worse than they really are"
$ws.Range("K6").Value = "You can use the directly measured ACM ratios of the two cohorts to establish the frailty at that age. That's the easiest."

# Update the active cell selection to match the saved view state (K6 -> L13).
$ws.Range("L13").Select()
